{"js": "// Update the date heading and every arithmetic-problem cell in the table.\n// Every <w:t> text run in the document changes; the surrounding\n// paragraph/run formatting (fonts, size, alignment) must stay untouched,\n// so we replace text in-place rather than clearing/re-creating runs.\n\nconst newDate = \"2024-01-12 Friday\";\n\n// New table contents, row-major (20 rows x 5 columns) - mirrors the\n// table.values shape exposed by the Word JS API.\nconst newValues = [\n  [\"73-73=\", \"90-29=\", \"3+6=\", \"17+56=\", \"45+38=\"],\n  [\"5+79=\", \"4+2=\", \"78-0=\", \"63+5=\", \"49+41=\"],\n  [\"44+34=\", \"55-5=\", \"77+15=\", \"2+62=\", \"89-49=\"],\n  [\"65-51=\", \"79-46=\", \"55-47=\", \"54+4=\", \"63-58=\"],\n  [\"37+42=\", \"2+39=\", \"16-12=\", \"70+10=\", \"22+63=\"],\n  [\"93-18=\", \"10+20=\", \"39-34=\", \"54+40=\", \"92-45=\"],\n  [\"54-28=\", \"33+22=\", \"76+15=\", \"96-90=\", \"0+85=\"],\n  [\"71+27=\", \"67+7=\", \"95+4=\", \"76-38=\", \"40+29=\"],\n  [\"45+24=\", \"28+35=\", \"69-21=\", \"50+36=\", \"26-26=\"],\n  [\"9+82=\", \"14+3=\", \"68+22=\", \"34-6=\", \"35+59=\"],\n  [\"18+17=\", \"84-41=\", \"88+1=\", \"97-87=\", \"69+8=\"],\n  [\"2+25=\", \"11+77=\", \"95-61=\", \"78+9=\", \"76-49=\"],\n  [\"93-88=\", \"65-26=\", \"86-46=\", \"4+61=\", \"3+85=\"],\n  [\"85-1=\", \"55-43=\", \"75-50=\", \"86+12=\", \"20-11=\"],\n  [\"33+7=\", \"78-71=\", \"31+5=\", \"37+56=\", \"37+9=\"],\n  [\"26+9=\", \"29+50=\", \"72-11=\", \"95-6=\", \"35+12=\"],\n  [\"67-60=\", \"59+18=\", \"86-58=\", \"57+32=\", \"47+4=\"],\n  [\"49-33=\", \"64-42=\", \"50+49=\", \"95-26=\", \"8+24=\"],\n  [\"94-48=\", \"32+67=\", \"48+41=\", \"57-26=\", \"19+39=\"],\n  [\"63+2=\", \"98-67=\", \"99-45=\", \"49+20=\", \"18+1=\"],\n];\n\n// 1) Update the date paragraph (first paragraph in the body), keeping its\n//    run/paragraph formatting intact.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titleParagraph = paragraphs.items[0];\ntitleParagraph.getRange().insertText(newDate, Word.InsertLocation.replace);\n\n// 2) Update every cell of the (single) table with the new values. Setting\n//    `.values` rewrites only the text runs inside each cell paragraph and\n//    preserves existing run/paragraph formatting.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = newValues;\n\nawait context.sync();\n", "ps1": "# Update the date heading and every arithmetic-problem cell in the table.\n# Every piece of visible text in the document changes, but the\n# paragraph/run formatting (fonts, size, alignment) must be preserved, so\n# we overwrite Range.Text in place instead of deleting/recreating runs.\n\n$d = $word.ActiveDocument\n\n$newDate = \"2024-01-12 Friday\"\n\n# New table contents, row-major (20 rows x 5 columns).\n$newValues = @(\n    @(\"73-73=\", \"90-29=\", \"3+6=\", \"17+56=\", \"45+38=\"),\n    @(\"5+79=\", \"4+2=\", \"78-0=\", \"63+5=\", \"49+41=\"),\n    @(\"44+34=\", \"55-5=\", \"77+15=\", \"2+62=\", \"89-49=\"),\n    @(\"65-51=\", \"79-46=\", \"55-47=\", \"54+4=\", \"63-58=\"),\n    @(\"37+42=\", \"2+39=\", \"16-12=\", \"70+10=\", \"22+63=\"),\n    @(\"93-18=\", \"10+20=\", \"39-34=\", \"54+40=\", \"92-45=\"),\n    @(\"54-28=\", \"33+22=\", \"76+15=\", \"96-90=\", \"0+85=\"),\n    @(\"71+27=\", \"67+7=\", \"95+4=\", \"76-38=\", \"40+29=\"),\n    @(\"45+24=\", \"28+35=\", \"69-21=\", \"50+36=\", \"26-26=\"),\n    @(\"9+82=\", \"14+3=\", \"68+22=\", \"34-6=\", \"35+59=\"),\n    @(\"18+17=\", \"84-41=\", \"88+1=\", \"97-87=\", \"69+8=\"),\n    @(\"2+25=\", \"11+77=\", \"95-61=\", \"78+9=\", \"76-49=\"),\n    @(\"93-88=\", \"65-26=\", \"86-46=\", \"4+61=\", \"3+85=\"),\n    @(\"85-1=\", \"55-43=\", \"75-50=\", \"86+12=\", \"20-11=\"),\n    @(\"33+7=\", \"78-71=\", \"31+5=\", \"37+56=\", \"37+9=\"),\n    @(\"26+9=\", \"29+50=\", \"72-11=\", \"95-6=\", \"35+12=\"),\n    @(\"67-60=\", \"59+18=\", \"86-58=\", \"57+32=\", \"47+4=\"),\n    @(\"49-33=\", \"64-42=\", \"50+49=\", \"95-26=\", \"8+24=\"),\n    @(\"94-48=\", \"32+67=\", \"48+41=\", \"57-26=\", \"19+39=\"),\n    @(\"63+2=\", \"98-67=\", \"99-45=\", \"49+20=\", \"18+1=\")\n)\n\n# 1) Update the date paragraph (first paragraph in the body), keeping its\n#    run/paragraph formatting intact.\n$d.Paragraphs.Item(1).Range.Text = $newDate\n\n# 2) Update every cell of the (single) table with the new values.\n$t = $d.Tables.Item(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$r - 1][$c - 1]\n    }\n}\n"}
